$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new observation record was added to the dataset, ahead of the existing
# rows. Insert a fresh row at position 2, which pushes the current rows
# 2 and 3 down to 3 and 4 respectively (their data is preserved as-is).
$ws.Rows("2:2").Insert()

# --- Populate the newly inserted row 2 --------------------------------
$ws.Cells.Item(2,1).Value  = 98838937          # A2  Id
$ws.Cells.Item(2,2).Value  = 104838            # B2  Taxonsorteringsordning
$ws.Cells.Item(2,3).Value  = "Ovaliderad"      # C2  Valideringsstatus
$ws.Cells.Item(2,4).Value  = "VU"              # D2  Rödlistade
$ws.Cells.Item(2,5).Value  = 219955            # E2  TaxonId
$ws.Cells.Item(2,6).Value  = "Slåttergubbe"    # F2  Artnamn
$ws.Cells.Item(2,7).Value  = "Arnica montana"  # G2  Vetenskapligt namn
$ws.Cells.Item(2,8).Value  = "L."              # H2  Auktor

# Antal ("100") must stay textual, not become the number 100.
$ws.Cells.Item(2,9).NumberFormat = "@"
$ws.Cells.Item(2,9).Value  = "100"             # I2  Antal

$ws.Cells.Item(2,10).Value = "plantor/tuvor"   # J2  Enhet

$ws.Cells.Item(2,16).Value = "Bondhyttan, Dlr"            # P2  Lokalnamn
$ws.Cells.Item(2,17).Value = 529961.0107277337             # Q2  Ost
$ws.Cells.Item(2,18).Value = 6685229.121118558              # R2  Nord
$ws.Cells.Item(2,19).Value = 10                             # S2  Noggrannhet
$ws.Cells.Item(2,20).Value = "Dalarna"                       # T2  Län
$ws.Cells.Item(2,21).Value = "Säter"                         # U2  Kommun
$ws.Cells.Item(2,22).Value = "Dalarna"                       # V2  Provins
$ws.Cells.Item(2,23).Value = "Silvberg"                      # W2  Församling
$ws.Cells.Item(2,24).Value = "W-Sät-0040"                    # X2  Externid

# Start-/slutdatum ("1989-01-01") must stay textual, not become a date serial.
$ws.Cells.Item(2,25).NumberFormat = "@"
$ws.Cells.Item(2,25).Value = "1989-01-01"      # Y2  Startdatum
$ws.Cells.Item(2,26).Value = "00:00"           # Z2  Starttid
$ws.Cells.Item(2,27).NumberFormat = "@"
$ws.Cells.Item(2,27).Value = "1989-01-01"      # AA2 Slutdatum
$ws.Cells.Item(2,28).Value = "00:00"           # AB2 Sluttid

$ws.Cells.Item(2,30).Value = $false            # AD2 Ej återfunnen
$ws.Cells.Item(2,31).Value = $false            # AE2 Osäker artbestämning
$ws.Cells.Item(2,33).Value = $false            # AG2 Ospontan

$ws.Cells.Item(2,49).Value = "Lennart Bratt"            # AW2 Rapportör
$ws.Cells.Item(2,50).Value = "Sören Nyström"             # AX2 Observatörer
$ws.Cells.Item(2,51).Value = "Floraväkteri Sverige"      # AY2 Projektnamn
